# The workbook holds curated metadata describing the "municipio-nombre"
# (column F) alongside "provincia-nombre" (column G) and "comarca-nombre"
# (column J). Re-processing the data with the newly curated dimensions
# reclassifies "municipio-nombre" from a measure to a proper SDMX
# dimension, matching the pattern already used by the sibling columns:
#   row 2 (concept type) : iaest-measure:municipio-nombre -> sdmx-dimension:refArea
#   row 3 (medida/dim)   : medida                          -> dim
#   row 4 (datatype/URI) : xsd:int                          -> URI-Municipio
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = "sdmx-dimension:refArea"
$ws.Range("F3").Value = "dim"
$ws.Range("F4").Value = "URI-Municipio"
